$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reference voltage (H1) used by the ADC-count formulas in
# column G (=ROUND(F/$H$1*1023,0)); this single input change cascades
# through the existing formulas and recalculates the dependent cells.
$ws.Range("H1").Value = 4.955

# Reset the view: scroll back to the top of the sheet and move the
# active selection to H2 (next to the cell that was just edited).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H2").Select()
